$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: 2023-07-04, with "y" marker for Nereocystis_luetkeana (column E)
$ws.Range("A13").Value = 2023
$ws.Range("B13").Value = 7
$ws.Range("C13").Value = 4
$ws.Range("E13").Value = "y"

# Row 14: 2023-08-01, with "y" markers for Saccharina_latissima (D) and Nereocystis_luetkeana (E)
$ws.Range("A14").Value = 2023
$ws.Range("B14").Value = 8
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = "y"
$ws.Range("E14").Value = "y"

# Update selection to match the final state (active cell C14)
$ws.Range("C14").Select()
